$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "Day"
$ws.Range("C1").Value = "Exchange"
$ws.Range("D1").Value = "Holiday"
$ws.Range("E1").Value = "Description"
$ws.Range("G11").Select()
